function Set-TextCell {
    param($ws, $addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextCell $ws "D2" "26.143.17"
Set-TextCell $ws "D3" "1.654.85"
Set-TextCell $ws "E3" "  +0.24%  "
Set-TextCell $ws "D5" "218.82"
Set-TextCell $ws "D6" "0.5239"
Set-TextCell $ws "E6" "  +0.19%  "
Set-TextCell $ws "D8" "0.2659"
Set-TextCell $ws "E9" "  +1.15%  "
Set-TextCell $ws "E10" "  +0.12%  "
Set-TextCell $ws "D11" "0.07692"
Set-TextCell $ws "E11" "  -1.31%  "
Set-TextCell $ws "D12" "4.638"
Set-TextCell $ws "E12" "  +3.55%  "
Set-TextCell $ws "D13" "1.653.31"
Set-TextCell $ws "E13" "  +0.44%  "
Set-TextCell $ws "D14" "1.882.65"
Set-TextCell $ws "E14" "  +0.26%  "
Set-TextCell $ws "D15" "0.5615"
Set-TextCell $ws "E15" "  +1.57%  "
Set-TextCell $ws "D16" "0.0₅8198"
Set-TextCell $ws "E16" "  +2.40%  "
Set-TextCell $ws "D17" "65.56"
Set-TextCell $ws "E17" "  +1.30%  "
Set-TextCell $ws "D18" "26.134.53"
Set-TextCell $ws "E18" "  +0.16%  "
Set-TextCell $ws "E19" "  -0.23%  "
Set-TextCell $ws "D20" "4.655"
Set-TextCell $ws "E20" "  +0.78%  "
Set-TextCell $ws "D21" "10.48"
Set-TextCell $ws "E21" "  +4.14%  "
Set-TextCell $ws "D22" "192.42"
Set-TextCell $ws "E22" "  -0.97%  "
Set-TextCell $ws "D23" "5.960"
Set-TextCell $ws "E23" "  +0.26%  "
Set-TextCell $ws "E24" "  -0.33%  "
Set-TextCell $ws "D25" "145.32"
Set-TextCell $ws "E25" "  -1.24%  "
Set-TextCell $ws "E26" "  -0.46%  "
Set-TextCell $ws "D27" "7.264"
Set-TextCell $ws "E27" "  +1.51%  "
Set-TextCell $ws "D28" "15.97"
Set-TextCell $ws "E28" "  +0.58%  "
Set-TextCell $ws "E29" "  +2.51%  "
Set-TextCell $ws "E30" "  -3.93%  "
Set-TextCell $ws "D31" "1.272"
Set-TextCell $ws "E31" "  +0.30%  "
Set-TextCell $ws "D32" "3.462"
Set-TextCell $ws "E32" "  -0.39%  "
Set-TextCell $ws "D33" "3.379"
Set-TextCell $ws "E33" "  +1.24%  "
Set-TextCell $ws "E34" "  -1.63%  "
Set-TextCell $ws "D35" "0.9545"
Set-TextCell $ws "E35" "  +0.63%  "
Set-TextCell $ws "E36" "  -0.59%  "
Set-TextCell $ws "D37" "2.402"
Set-TextCell $ws "E37" "  -0.46%  "
Set-TextCell $ws "D38" "0.5663"
Set-TextCell $ws "E38" "  +0.11%  "
Set-TextCell $ws "D39" "0.01582"
Set-TextCell $ws "E39" "  -0.45%  "
Set-TextCell $ws "D40" "5.870"
Set-TextCell $ws "E40" "  -0.93%  "
Set-TextCell $ws "E41" "  -0.27%  "
Set-TextCell $ws "D42" "0.8364"
Set-TextCell $ws "E42" "  -0.65%  "
Set-TextCell $ws "E43" "  -3.02%  "
Set-TextCell $ws "D44" "101.26"
Set-TextCell $ws "E44" "  -2.02%  "
Set-TextCell $ws "D45" "1.793.61"
Set-TextCell $ws "E45" "  +0.28%  "
Set-TextCell $ws "D46" "57.88"
Set-TextCell $ws "E46" "  +0.88%  "
Set-TextCell $ws "E47" "  +2.97%  "
Set-TextCell $ws "D48" "0.9999"
Set-TextCell $ws "E48" "  -0.89%  "
Set-TextCell $ws "B49" "Mantle"
Set-TextCell $ws "C49" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell $ws "D49" "0.4339"
Set-TextCell $ws "E49" "  -1.30%  "
Set-TextCell $ws "B50" "EnergySwap"
Set-TextCell $ws "C50" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws "D50" "7.994"
Set-TextCell $ws "E50" "  +0.73%  "
Set-TextCell $ws "D51" "0.05193"
Set-TextCell $ws "E51" "  -3.81%  "
